$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 1 header: shift existing J1/K1/L1 -> K1/L1/M1 ---
$ws.Range("J1").ClearContents()
$ws.Range("K1").ClearContents()
$ws.Range("L1").ClearContents()

$ws.Range("K1").Value = "zawor1"
$ws.Range("L1").Value = "zawor2"
$ws.Range("M1").Value = "zawor3"

# New text values - write them in the same order they were first entered so
# the shared-string table indices come out in the same order as the target.
$ws.Range("I2").Value = "52s"
$ws.Range("I1").Value = "czas do 5 cm"
$ws.Range("T1").Value = "wymiary zbiornika"
$ws.Range("U3").Value = "w 30: 41.5"
$ws.Range("V3").Value = "w 5: 15.5"
$ws.Range("W3").Value = "głębokość: 5"
$ws.Range("A10").Value = "objętość"

# --- Row 2 (remaining numeric / formula cells) ---
$ws.Range("H2").Value = 8.69
$ws.Range("J2").Value = 1
$ws.Range("K2").Formula = "=B10/144.8"
$ws.Range("L2").Formula = "=C10/119.4"
$ws.Range("O2").Formula = "=303.7-52"
$ws.Range("T2").Value = 1
$ws.Range("U2").Value = 31
$ws.Range("V2").Value = 5
$ws.Range("W2").Value = 40

# --- Row 3 ---
$ws.Range("H3").Formula = "=`$B`$10/O3"
$ws.Range("I3").Value = 19.08
$ws.Range("J3").Value = 0.9
$ws.Range("K3").Formula = "=B10/150.6"
$ws.Range("L3").Formula = "=C10/124.43"
$ws.Range("O3").Formula = "=116.2-I3"
$ws.Range("T3").Value = 2
$ws.Range("X3").Formula = "=(41.5+10)/2*5*30"
$ws.Range("Y3").Formula = "=(15.5+10)/2*5*5"

# --- Row 4 ---
$ws.Range("H4").Formula = "=`$B`$10/O4"
$ws.Range("I4").Value = 13.5
$ws.Range("J4").Value = 0.8
$ws.Range("K4").Formula = "=B10/157"
$ws.Range("O4").Formula = "=78.6-I4"

# --- Row 5 ---
$ws.Range("H5").Formula = "=`$B`$10/O5"
$ws.Range("I5").Value = 6.79
$ws.Range("J5").Value = 0.7
$ws.Range("K5").Formula = "=B10/177.6"
$ws.Range("O5").Formula = "=57.8-I5"

# --- Row 6 ---
$ws.Range("H6").Formula = "=`$B`$10/O6"
$ws.Range("I6").Value = 5.77
$ws.Range("J6").Value = 0.6
$ws.Range("K6").Formula = "=B10/214.5"
$ws.Range("O6").Formula = "=48.83-I6"

# --- Row 7 ---
$ws.Range("H7").Formula = "=`$B`$10/O7"
$ws.Range("I7").Value = 6.39
$ws.Range("J7").Value = 0.5
$ws.Range("K7").Formula = "=B10/265.4"
$ws.Range("O7").Formula = "=44.21-I7"

# --- Row 8 ---
$ws.Range("H8").Formula = "=`$B`$10/O8"
$ws.Range("I8").Value = 5.42
$ws.Range("J8").Value = 0.4
$ws.Range("K8").Formula = "=B10/478.9"
$ws.Range("O8").Formula = "=39.1-I8"

# --- Row 9 ---
$ws.Range("H9").Formula = "=`$B`$10/O9"
$ws.Range("I9").Value = 4.49
$ws.Range("J9").Value = 0.3
$ws.Range("O9").Formula = "=34.7-I9"

# --- Row 10 ---
$ws.Range("B10").Formula = "=25*31*5"
$ws.Range("C10").Formula = "=X3-Y3"
$ws.Range("H10").Formula = "=`$B`$10/O10"
$ws.Range("I10").Value = 5.47
$ws.Range("O10").Formula = "=33.19-I10"

# --- Column widths (match Excel's "best fit" auto-sized widths) ---
$ws.Columns.Item(1).ColumnWidth = 12.0
$ws.Columns.Item(9).ColumnWidth = 10.83
$ws.Columns.Item(23).ColumnWidth = 14.0

# --- View settings ---
$excel.ActiveWindow.Zoom = 110
$ws.Range("A11").Select()
